# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# header row / data, mirroring the formatting already used by column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells (rows 2-3)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

# Copy the header formatting (bold font, borders, centered) from H1 onto
# the two new header cells so I1/J1 match the style used by the rest of
# the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
